$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "97.837.21"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.95%  "

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.409.04"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +3.18%  "

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "256.16"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.36%  "

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "658.91"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +5.53%  "

$ws.Range("E7").Value = "  +0.24%  "

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.432"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +7.96%  "

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "1.06"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +9.25%  "

$ws.Range("E10").Value = "  -0.06%  "

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "3.405.37"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +3.08%  "

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.213"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +6.66%  "

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "41.87"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +5.41%  "

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "6.29"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +15.11%  "

$ws.Range("E15").Value = "  +4.71%  "

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "97.397.77"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -1.05%  "

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "4.027.18"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +2.93%  "

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "8.56"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +35.57%  "

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "3.395.51"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +2.73%  "

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "17.41"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +11.57%  "

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "0.497"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +64.25%  "

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "10.84"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +14.59%  "

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "3.46"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.06%  "

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "509.98"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +4.78%  "

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.0000206"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +2.36%  "

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "6.16"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +9.87%  "

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "96.94"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +9.50%  "

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "12.70"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +5.71%  "

$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.151"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +11.55%  "

$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "11.56"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +14.98%  "

$ws.Range("B31").Value = "Dai"
$ws.Range("C31").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.995"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.42%  "

$ws.Range("B32").Value = "Cronos"
$ws.Range("C32").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.193"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +2.57%  "

$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.993"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.35%  "

$ws.Range("B34").Value = "PolygonEcosystemToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.565"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +21.46%  "

$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "29.75"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +6.13%  "

$ws.Range("B36").Value = "PancakeSwap"
$ws.Range("C36").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "2.18"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +11.96%  "

$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "7.81"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +8.81%  "

$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.157"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +6.71%  "

$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "512.17"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +4.78%  "

$ws.Range("B40").Value = "Fetch.AI"
$ws.Range("C40").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "1.39"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +12.28%  "

$ws.Range("B41").Value = "WhiteBITCoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "24.70"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.58%  "

$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.0432"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +30.40%  "

$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.852"
$c.Style = "Normal"
$ws.Range("E43").Value = "  +8.25%  "

$ws.Range("B44").Value = "MantraDAO"
$ws.Range("C44").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "3.68"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +1.06%  "

$ws.Range("B45").Value = "dogwifhat"
$ws.Range("C45").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "3.31"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +6.88%  "

$ws.Range("B46").Value = "Cosmos"
$ws.Range("C46").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "8.30"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +13.25%  "

$ws.Range("B47").Value = "USDe"
$ws.Range("C47").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -0.03%  "

$ws.Range("B48").Value = "Filecoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "5.32"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +12.58%  "

$ws.Range("B49").Value = "ImmutableX"
$ws.Range("C49").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.57"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +15.87%  "

$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "2.08"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +7.49%  "

$ws.Range("B51").Value = "OKB"
$ws.Range("C51").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "50.46"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +10.43%  "
